$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data range A37:B44 contains two rows (angle=276 and angle=350) that were
# out of ascending order relative to the rest of the "ty (deg)" angle column.
# Re-write the range A37:B44 in ascending angle order.
$values = @(
    @(276, 3.6),
    @(280, 3.11),
    @(285, 2.83),
    @(290, 2.57),
    @(295, 2.21),
    @(300, 1.89),
    @(306, 1.61),
    @(350, 0.33)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 37 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Update the active cell/selection to D41, matching the new cursor position
# left after performing the edit.
$ws.Range("D41").Select()
